# Applies the "Natmi following Dr Hou advice" update to the LR-pairs sheet.
# Rewrites the 16-row x 20-column data table (rows 2-17) so every
# Sending-cluster x Target-cluster combination (including the previously
# missing "Target cluster = ECs" rows) is present with refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 16,20

# Row 2: ECs -> ECs
$data[0,0] = "ECs"
$data[0,1] = "Col18a1"
$data[0,2] = "Gpc1"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 11.122774
$data[0,7] = 33.368322
$data[0,8] = 0.2449652610853511
$data[0,9] = 0.2449652610853511
$data[0,10] = 2
$data[0,11] = 0.6666666666666666
$data[0,12] = 1.628421
$data[0,13] = 4.885263
$data[0,14] = 0.048329411442081
$data[0,15] = 0.048329411442081
$data[0,16] = 18.112558759854
$data[0,17] = 163.013028838686
$data[0,18] = 0.01183902689201073
$data[0,19] = 0.01183902689201073

# Row 3: ECs -> FAPs
$data[1,0] = "ECs"
$data[1,1] = "Col18a1"
$data[1,2] = "Gpc1"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 11.122774
$data[1,7] = 33.368322
$data[1,8] = 0.2449652610853511
$data[1,9] = 0.2449652610853511
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 9.459065000000001
$data[1,13] = 28.377195
$data[1,14] = 0.2807327123897247
$data[1,15] = 0.2807327123897247
$data[1,16] = 105.21104224631
$data[1,17] = 946.89938021679
$data[1,18] = 0.06876976218574771
$data[1,19] = 0.06876976218574769

# Row 4: ECs -> M2
$data[2,0] = "ECs"
$data[2,1] = "Col18a1"
$data[2,2] = "Gpc1"
$data[2,3] = "M2"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 11.122774
$data[2,7] = 33.368322
$data[2,8] = 0.2449652610853511
$data[2,9] = 0.2449652610853511
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 0.6418243333333334
$data[2,13] = 1.925473
$data[2,14] = 0.01904850912583786
$data[2,15] = 0.01904850912583786
$data[2,16] = 7.138867007367335
$data[2,17] = 64.249803066306
$data[2,18] = 0.004666223011297564
$data[2,19] = 0.004666223011297564

# Row 5: ECs -> sCs
$data[3,0] = "ECs"
$data[3,1] = "Col18a1"
$data[3,2] = "Gpc1"
$data[3,3] = "sCs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 11.122774
$data[3,7] = 33.368322
$data[3,8] = 0.2449652610853511
$data[3,9] = 0.2449652610853511
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 21.96489266666667
$data[3,13] = 65.894678
$data[3,14] = 0.6518893670423564
$data[3,15] = 0.6518893670423563
$data[3,16] = 244.3105370655907
$data[3,17] = 2198.794833590316
$data[3,18] = 0.1596902489962951
$data[3,19] = 0.1596902489962951

# Row 6: FAPs -> ECs
$data[4,0] = "FAPs"
$data[4,1] = "Col18a1"
$data[4,2] = "Gpc1"
$data[4,3] = "ECs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 18.220714
$data[4,7] = 54.662142
$data[4,8] = 0.4012885600454987
$data[4,9] = 0.4012885600454988
$data[4,10] = 2
$data[4,11] = 0.6666666666666666
$data[4,12] = 1.628421
$data[4,13] = 4.885263
$data[4,14] = 0.048329411442081
$data[4,15] = 0.048329411442081
$data[4,16] = 29.670993312594
$data[4,17] = 267.038939813346
$data[4,18] = 0.01939403992543913
$data[4,19] = 0.01939403992543914

# Row 7: FAPs -> FAPs
$data[5,0] = "FAPs"
$data[5,1] = "Col18a1"
$data[5,2] = "Gpc1"
$data[5,3] = "FAPs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 18.220714
$data[5,7] = 54.662142
$data[5,8] = 0.4012885600454987
$data[5,9] = 0.4012885600454988
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 9.459065000000001
$data[5,13] = 28.377195
$data[5,14] = 0.2807327123897247
$data[5,15] = 0.2807327123897247
$data[5,16] = 172.35091807241
$data[5,17] = 1551.15826265169
$data[5,18] = 0.1126548259125398
$data[5,19] = 0.1126548259125398

# Row 8: FAPs -> M2
$data[6,0] = "FAPs"
$data[6,1] = "Col18a1"
$data[6,2] = "Gpc1"
$data[6,3] = "M2"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 18.220714
$data[6,7] = 54.662142
$data[6,8] = 0.4012885600454987
$data[6,9] = 0.4012885600454988
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 0.6418243333333334
$data[6,13] = 1.925473
$data[6,14] = 0.01904850912583786
$data[6,15] = 0.01904850912583786
$data[6,16] = 11.69449761590734
$data[6,17] = 105.250478543166
$data[6,18] = 0.007643948798121017
$data[6,19] = 0.007643948798121017

# Row 9: FAPs -> sCs
$data[7,0] = "FAPs"
$data[7,1] = "Col18a1"
$data[7,2] = "Gpc1"
$data[7,3] = "sCs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 18.220714
$data[7,7] = 54.662142
$data[7,8] = 0.4012885600454987
$data[7,9] = 0.4012885600454988
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 21.96489266666667
$data[7,13] = 65.894678
$data[7,14] = 0.6518893670423564
$data[7,15] = 0.6518893670423563
$data[7,16] = 400.2160273200307
$data[7,17] = 3601.944245880276
$data[7,18] = 0.2615957454093988
$data[7,19] = 0.2615957454093988

# Row 10: M2 -> ECs
$data[8,0] = "M2"
$data[8,1] = "Col18a1"
$data[8,2] = "Gpc1"
$data[8,3] = "ECs"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 0.1189986666666667
$data[8,7] = 0.356996
$data[8,8] = 0.002620797603979787
$data[8,9] = 0.002620797603979787
$data[8,10] = 2
$data[8,11] = 0.6666666666666666
$data[8,12] = 1.628421
$data[8,13] = 4.885263
$data[8,14] = 0.048329411442081
$data[8,15] = 0.048329411442081
$data[8,16] = 0.193779927772
$data[8,17] = 1.744019349948
$data[8,18] = 0.0001266616057091592
$data[8,19] = 0.0001266616057091592

# Row 11: M2 -> FAPs
$data[9,0] = "M2"
$data[9,1] = "Col18a1"
$data[9,2] = "Gpc1"
$data[9,3] = "FAPs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 0.1189986666666667
$data[9,7] = 0.356996
$data[9,8] = 0.002620797603979787
$data[9,9] = 0.002620797603979787
$data[9,10] = 3
$data[9,11] = 1
$data[9,12] = 9.459065000000001
$data[9,13] = 28.377195
$data[9,14] = 0.2807327123897247
$data[9,15] = 0.2807327123897247
$data[9,16] = 1.125616122913333
$data[9,17] = 10.13054510622
$data[9,18] = 0.0007357436199897371
$data[9,19] = 0.000735743619989737

# Row 12: M2 -> M2
$data[10,0] = "M2"
$data[10,1] = "Col18a1"
$data[10,2] = "Gpc1"
$data[10,3] = "M2"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 0.1189986666666667
$data[10,7] = 0.356996
$data[10,8] = 0.002620797603979787
$data[10,9] = 0.002620797603979787
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 0.6418243333333334
$data[10,13] = 1.925473
$data[10,14] = 0.01904850912583786
$data[10,15] = 0.01904850912583786
$data[10,16] = 0.07637623990088889
$data[10,17] = 0.687386159108
$data[10,18] = 0.00004992228707638296
$data[10,19] = 0.00004992228707638296

# Row 13: M2 -> sCs
$data[11,0] = "M2"
$data[11,1] = "Col18a1"
$data[11,2] = "Gpc1"
$data[11,3] = "sCs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 0.1189986666666667
$data[11,7] = 0.356996
$data[11,8] = 0.002620797603979787
$data[11,9] = 0.002620797603979787
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 21.96489266666667
$data[11,13] = 65.894678
$data[11,14] = 0.6518893670423564
$data[11,15] = 0.6518893670423563
$data[11,16] = 2.613792940809778
$data[11,17] = 23.524136467288
$data[11,18] = 0.001708470091204507
$data[11,19] = 0.001708470091204507

# Row 14: sCs -> ECs
$data[12,0] = "sCs"
$data[12,1] = "Col18a1"
$data[12,2] = "Gpc1"
$data[12,3] = "ECs"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 15.943029
$data[12,7] = 47.829087
$data[12,8] = 0.3511253812651704
$data[12,9] = 0.3511253812651704
$data[12,10] = 2
$data[12,11] = 0.6666666666666666
$data[12,12] = 1.628421
$data[12,13] = 4.885263
$data[12,14] = 0.048329411442081
$data[12,15] = 0.048329411442081
$data[12,16] = 25.961963227209
$data[12,17] = 233.657669044881
$data[12,18] = 0.01696968301892198
$data[12,19] = 0.01696968301892198

# Row 15: sCs -> FAPs
$data[13,0] = "sCs"
$data[13,1] = "Col18a1"
$data[13,2] = "Gpc1"
$data[13,3] = "FAPs"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 15.943029
$data[13,7] = 47.829087
$data[13,8] = 0.3511253812651704
$data[13,9] = 0.3511253812651704
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 9.459065000000001
$data[13,13] = 28.377195
$data[13,14] = 0.2807327123897247
$data[13,15] = 0.2807327123897247
$data[13,16] = 150.806147607885
$data[13,17] = 1357.255328470965
$data[13,18] = 0.09857238067144752
$data[13,19] = 0.09857238067144751

# Row 16: sCs -> M2
$data[14,0] = "sCs"
$data[14,1] = "Col18a1"
$data[14,2] = "Gpc1"
$data[14,3] = "M2"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 15.943029
$data[14,7] = 47.829087
$data[14,8] = 0.3511253812651704
$data[14,9] = 0.3511253812651704
$data[14,10] = 3
$data[14,11] = 1
$data[14,12] = 0.6418243333333334
$data[14,13] = 1.925473
$data[14,14] = 0.01904850912583786
$data[14,15] = 0.01904850912583786
$data[14,16] = 10.232623959239
$data[14,17] = 92.09361563315102
$data[14,18] = 0.006688415029342896
$data[14,19] = 0.006688415029342896

# Row 17: sCs -> sCs
$data[15,0] = "sCs"
$data[15,1] = "Col18a1"
$data[15,2] = "Gpc1"
$data[15,3] = "sCs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 15.943029
$data[15,7] = 47.829087
$data[15,8] = 0.3511253812651704
$data[15,9] = 0.3511253812651704
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 21.96489266666667
$data[15,13] = 65.894678
$data[15,14] = 0.6518893670423564
$data[15,15] = 0.6518893670423563
$data[15,16] = 350.186920766554
$data[15,17] = 3151.682286898986
$data[15,18] = 0.228894902545458
$data[15,19] = 0.228894902545458

$ws.Range("A2:T17").Value = $data

